$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "GNG_TO-1650477808601557"
$wb.Worksheets.Item(2).Name = "NB_TO-16504778109205878"
$wb.Worksheets.Item(3).Name = "RS_TO-16504778109215536"
$wb.Worksheets.Item(4).Name = "TOL_TO-1650477810969557"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504778110305886"

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504778085725572.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778085845876.csv"
$ws1.Range("B4").Value = "go_stims-1650477808585553.csv"
$ws1.Range("B5").Value = "GNG_stims-1650477808600587.csv"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-1650477809599587.csv"
$ws2.Range("B3").Value = "ZB-match_7-16504778093165898.csv"
$ws2.Range("B4").Value = "ZB-match_3-16504778086145883.csv"
$ws2.Range("B5").Value = "OB-1650477809691594.csv"
$ws2.Range("B6").Value = "TB-16504778105525873.csv"
$ws2.Range("B7").Value = "TB-1650477810238591.csv"
$ws2.Range("B8").Value = "OB-16504778101585913.csv"
$ws2.Range("B9").Value = "ZB-match_8-16504778089675555.csv"
$ws2.Range("B10").Value = "TB-16504778109005868.csv"

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504778109365807.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778109235554.csv"
$ws4.Range("B4").Value = "MM_stims-1650477810952555.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778109375553.csv"
$ws4.Range("B6").Value = "MM_stims-16504778109685552.csv"
$ws4.Range("B7").Value = "ZM_stims-1650477810952555.csv"

$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16504778109715543.csv"
$ws5.Range("B3").Value = "vSAT_stims-16504778110155919.csv"
$ws5.Range("B4").Value = "vSAT_stims-16504778110005894.csv"
$ws5.Range("B5").Value = "SAT_stims-16504778109845574.csv"
